# close #125: Correct mandatory value columns assumptions in valores.xlsx
#
# The "1-2015" column (column C) was an extra/incorrect value column that
# should not have been part of the mandatory value columns - remove it
# entirely. Excel shifts every column to its right (D:K) one place to the
# left to fill the gap, so the sheet ends up spanning A1:J20 instead of
# A1:K20, and the headers/values that used to live in D:K now live in C:J.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(11).Select()
$ws.Columns.Item(3).Delete()
